# "just changing variable names for running" -- rework the percentile
# threshold checking sheet: rename/expand headers, add a "one JJA" vs
# "twenty years of JJAs" comparison block (columns I/J) and a new
# "No. of data points" pair of columns (F/G), plus new percentile rows
# (95, 97, 99.75, 50).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old table first so that stale shared-formula groupings
# from the original layout don't leak into the rebuilt one.
$ws.Range("A1:G10").ClearContents()

# ---- Row 1 : headers (order matters for shared-string indices) --------
$ws.Range("A1").Value = "Total data points"
$ws.Range("C1").Value = "Percentile"
$ws.Range("D1").Value = "Percentile fraction"
$ws.Range("E1").Value = "Remainder"

# ---- Column I / J helper values (labels referenced later by F1/G1) ----
$ws.Range("I2").Value = "Points in twenty years worth of JJAs"
$ws.Range("I3").Value = "Points in one JJA"

$ws.Range("F1").Value = "No. of data points in one JJA"
$ws.Range("G1").Value = "No. of data points in twenty years of JJAs"

# ---- Column A / J helper formulas --------------------------------------
$ws.Range("A2").Formula = '=1*90*24'

$ws.Range("J2").Formula = '=20*90*24'
$ws.Range("J3").Formula = '=24*90'
$ws.Range("J4").Formula = '=J2/J3'

# ---- Percentile column (C) --------------------------------------------
$ws.Range("C2").Value = 95
$ws.Range("C3").Value = 97
$ws.Range("C4").Value = 99
$ws.Range("C5").Value = 99.5
$ws.Range("C6").Value = 99.75
$ws.Range("C7").Value = 99.9
$ws.Range("C8").Value = 99.95
$ws.Range("C9").Value = 99.99
$ws.Range("C10").Value = 50

# ---- Percentile fraction (D) -------------------------------------------
# D2:D3 share one formula group, D4 is standalone, D5:D10 share another.
$ws.Range("D2:D3").Formula = '=C2/100'
$ws.Range("D4").Formula = '=C4/100'
$ws.Range("D5:D10").Formula = '=C5/100'

# ---- Remainder (E) ------------------------------------------------------
$ws.Range("E2:E3").Formula = '=1-D2'
$ws.Range("E4").Formula = '=1-D4'
$ws.Range("E5:E10").Formula = '=1-D5'

# ---- No. of data points in one JJA (F) -----------------------------------
$ws.Range("F2").Formula = '=E2*$A$2'
$ws.Range("F3").Formula = '=E3*$A$2'
$ws.Range("F4").Formula = '=E4*$A$2'
$ws.Range("F5").Formula = '=E5*$A$2'
$ws.Range("F6:F10").Formula = '=E6*$A$2'

# ---- No. of data points in twenty years of JJAs (G) ----------------------
$ws.Range("G2:G3").Formula = '=E2*$J$2'
$ws.Range("G4").Formula = '=E4*$J$2'
$ws.Range("G5:G10").Formula = '=E5*$J$2'

# ---- Column widths (approximate best-fit sizing from the original file) --
$ws.Columns("A:A").ColumnWidth = 6.28515625
$ws.Columns("F:F").ColumnWidth = 25.140625
$ws.Columns("G:G").ColumnWidth = 36.140625
$ws.Columns("I:I").ColumnWidth = 31.7109375

# ---- Selection matches the saved view in the edited workbook -------------
$ws.Range("G19").Select()
